$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 31 de Julio de 2020 a las 20:17"

$ws.Cells.Item(4, 2).Value = 4667579
$ws.Cells.Item(4, 3).Value = 32594
$ws.Cells.Item(4, 4).Value = 2289873
$ws.Cells.Item(4, 5).Value = 2221712
$ws.Cells.Item(4, 7).Value = 709
$ws.Cells.Item(4, 8).Value = 155994

$ws.Cells.Item(5, 2).Value = 2625612
$ws.Cells.Item(5, 3).Value = 11823
$ws.Cells.Item(5, 5).Value = 709910
$ws.Cells.Item(5, 7).Value = 230
$ws.Cells.Item(5, 8).Value = 91607

$ws.Cells.Item(6, 2).Value = 1694705
$ws.Cells.Item(6, 3).Value = 55355
$ws.Cells.Item(6, 4).Value = 1094279
$ws.Cells.Item(6, 5).Value = 563876
$ws.Cells.Item(6, 7).Value = 764
$ws.Cells.Item(6, 8).Value = 36550

$ws.Cells.Item(11, 2).Value = 355667
$ws.Cells.Item(11, 3).Value = 2131
$ws.Cells.Item(11, 4).Value = 328327
$ws.Cells.Item(11, 5).Value = 17883
$ws.Cells.Item(11, 7).Value = 80
$ws.Cells.Item(11, 8).Value = 9457

$ws.Cells.Item(18, 4).Value = 199974
$ws.Cells.Item(18, 5).Value = 12422

$ws.Cells.Item(20, 2).Value = 230873
$ws.Cells.Item(20, 3).Value = 982
$ws.Cells.Item(20, 4).Value = 214535
$ws.Cells.Item(20, 5).Value = 10647
$ws.Cells.Item(20, 7).Value = 17
$ws.Cells.Item(20, 8).Value = 5691

$ws.Cells.Item(22, 2).Value = 187919
$ws.Cells.Item(22, 3).Value = 1346
$ws.Cells.Item(22, 5).Value = 76154
$ws.Cells.Item(22, 7).Value = 11
$ws.Cells.Item(22, 8).Value = 30265

$ws.Cells.Item(25, 2).Value = 116116
$ws.Cells.Item(25, 3).Value = 317
$ws.Cells.Item(25, 5).Value = 6153
$ws.Cells.Item(25, 7).Value = 4
$ws.Cells.Item(25, 8).Value = 8933

$ws.Cells.Item(31, 2).Value = 85355
$ws.Cells.Item(31, 3).Value = 985
$ws.Cells.Item(31, 4).Value = 36044
$ws.Cells.Item(31, 5).Value = 43654

$ws.Cells.Item(33, 2).Value = 80422
$ws.Cells.Item(33, 3).Value = 45
$ws.Cells.Item(33, 7).Value = 9
$ws.Cells.Item(33, 8).Value = 5743

$ws.Cells.Item(36, 2).Value = 70970
$ws.Cells.Item(36, 3).Value = 934
$ws.Cells.Item(36, 4).Value = 43850
$ws.Cells.Item(36, 5).Value = 26608
$ws.Cells.Item(36, 7).Value = 12
$ws.Cells.Item(36, 8).Value = 512

$ws.Cells.Item(60, 2).Value = 30394
$ws.Cells.Item(60, 3).Value = 563
$ws.Cells.Item(60, 4).Value = 20537
$ws.Cells.Item(60, 5).Value = 8647
$ws.Cells.Item(60, 7).Value = 10
$ws.Cells.Item(60, 8).Value = 1210

$ws.Cells.Item(61, 2).Value = 26065
$ws.Cells.Item(61, 3).Value = 38
$ws.Cells.Item(61, 5).Value = 938

$ws.Cells.Item(64, 1).Value = "Marruecos"
$ws.Cells.Item(64, 2).Value = 24322
$ws.Cells.Item(64, 3).Value = 1063
$ws.Cells.Item(64, 4).Value = 17658
$ws.Cells.Item(64, 5).Value = 6311
$ws.Cells.Item(64, 7).Value = 7
$ws.Cells.Item(64, 8).Value = 353

$ws.Cells.Item(65, 1).Value = "Uzbekistan"
$ws.Cells.Item(65, 2).Value = 23773
$ws.Cells.Item(65, 3).Value = 502
$ws.Cells.Item(65, 4).Value = 14204
$ws.Cells.Item(65, 5).Value = 9430
$ws.Cells.Item(65, 7).Value = 3
$ws.Cells.Item(65, 8).Value = 139

$ws.Cells.Item(70, 1).Value = "Etiopia"
$ws.Cells.Item(70, 2).Value = 17530
$ws.Cells.Item(70, 3).Value = 915
$ws.Cells.Item(70, 4).Value = 6950
$ws.Cells.Item(70, 5).Value = 10306
$ws.Cells.Item(70, 7).Value = 11
$ws.Cells.Item(70, 8).Value = 274

$ws.Cells.Item(71, 1).Value = "Costa Rica"
$ws.Cells.Item(71, 2).Value = 17290
$ws.Cells.Item(71, 4).Value = 4280
$ws.Cells.Item(71, 5).Value = 12870
$ws.Cells.Item(71, 8).Value = 140

$ws.Cells.Item(72, 1).Value = "Camerun"
$ws.Cells.Item(72, 2).Value = 17255
$ws.Cells.Item(72, 3).Value = 0
$ws.Cells.Item(72, 4).Value = 15320
$ws.Cells.Item(72, 5).Value = 1544
$ws.Cells.Item(72, 7).Value = 0
$ws.Cells.Item(72, 8).Value = 391

$ws.Cells.Item(73, 1).Value = "Australia"
$ws.Cells.Item(73, 2).Value = 16905
$ws.Cells.Item(73, 3).Value = 602
$ws.Cells.Item(73, 4).Value = 9982
$ws.Cells.Item(73, 5).Value = 6726
$ws.Cells.Item(73, 7).Value = 8
$ws.Cells.Item(73, 8).Value = 197

$ws.Cells.Item(74, 1).Value = "El Salvador"
$ws.Cells.Item(74, 2).Value = 16632
$ws.Cells.Item(74, 3).Value = 402
$ws.Cells.Item(74, 4).Value = 8362
$ws.Cells.Item(74, 5).Value = 7822
$ws.Cells.Item(74, 7).Value = 9
$ws.Cells.Item(74, 8).Value = 448

$ws.Cells.Item(75, 2).Value = 16474
$ws.Cells.Item(75, 3).Value = 103
$ws.Cells.Item(75, 4).Value = 11559
$ws.Cells.Item(75, 5).Value = 4533
$ws.Cells.Item(75, 7).Value = 3
$ws.Cells.Item(75, 8).Value = 382

$ws.Cells.Item(95, 2).Value = 6695
$ws.Cells.Item(95, 3).Value = 79
$ws.Cells.Item(95, 4).Value = 5192
$ws.Cells.Item(95, 5).Value = 1389

$ws.Cells.Item(98, 1).Value = "Albania"
$ws.Cells.Item(98, 2).Value = 5276
$ws.Cells.Item(98, 3).Value = 79
$ws.Cells.Item(98, 4).Value = 2952
$ws.Cells.Item(98, 5).Value = 2167
$ws.Cells.Item(98, 7).Value = 3
$ws.Cells.Item(98, 8).Value = 157

$ws.Cells.Item(99, 1).Value = "Paraguay"
$ws.Cells.Item(99, 2).Value = 5207
$ws.Cells.Item(99, 4).Value = 3250
$ws.Cells.Item(99, 5).Value = 1910
$ws.Cells.Item(99, 8).Value = 47

$ws.Cells.Item(103, 1).Value = "Libano"
$ws.Cells.Item(103, 2).Value = 4555
$ws.Cells.Item(103, 3).Value = 221
$ws.Cells.Item(103, 4).Value = 1761
$ws.Cells.Item(103, 5).Value = 2733
$ws.Cells.Item(103, 7).Value = 4
$ws.Cells.Item(103, 8).Value = 61

$ws.Cells.Item(104, 1).Value = "Hungria"
$ws.Cells.Item(104, 2).Value = 4505
$ws.Cells.Item(104, 3).Value = 21
$ws.Cells.Item(104, 4).Value = 3353
$ws.Cells.Item(104, 5).Value = 556
$ws.Cells.Item(104, 8).Value = 596

$ws.Cells.Item(105, 1).Value = "Grecia"
$ws.Cells.Item(105, 2).Value = 4477
$ws.Cells.Item(105, 3).Value = 76
$ws.Cells.Item(105, 4).Value = 1374
$ws.Cells.Item(105, 5).Value = 2897
$ws.Cells.Item(105, 7).Value = 3
$ws.Cells.Item(105, 8).Value = 206

$ws.Cells.Item(107, 2).Value = 3793
$ws.Cells.Item(107, 3).Value = 74
$ws.Cells.Item(107, 4).Value = 2607
$ws.Cells.Item(107, 5).Value = 1170

$ws.Cells.Item(118, 2).Value = 2815
$ws.Cells.Item(118, 3).Value = 1
$ws.Cells.Item(118, 5).Value = 413

$ws.Cells.Item(122, 2).Value = 2451
$ws.Cells.Item(122, 3).Value = 33
$ws.Cells.Item(122, 4).Value = 1756
$ws.Cells.Item(122, 5).Value = 672

$ws.Cells.Item(132, 1).Value = "Mozambique"
$ws.Cells.Item(132, 2).Value = 1864
$ws.Cells.Item(132, 3).Value = 56
$ws.Cells.Item(132, 4).Value = 641
$ws.Cells.Item(132, 5).Value = 1212
$ws.Cells.Item(132, 8).Value = 11

$ws.Cells.Item(133, 1).Value = "Sierra Leona"
$ws.Cells.Item(133, 2).Value = 1823
$ws.Cells.Item(133, 3).Value = 5
$ws.Cells.Item(133, 4).Value = 1362
$ws.Cells.Item(133, 5).Value = 394
$ws.Cells.Item(133, 8).Value = 67

$ws.Cells.Item(135, 2).Value = 1728
$ws.Cells.Item(135, 3).Value = 2
$ws.Cells.Item(135, 4).Value = 862
$ws.Cells.Item(135, 5).Value = 373
$ws.Cells.Item(135, 7).Value = 6
$ws.Cells.Item(135, 8).Value = 493

$ws.Cells.Item(136, 2).Value = 1650
$ws.Cells.Item(136, 3).Value = 43
$ws.Cells.Item(136, 4).Value = 1080
$ws.Cells.Item(136, 5).Value = 544

$ws.Cells.Item(144, 2).Value = 1154
$ws.Cells.Item(144, 3).Value = 7
$ws.Cells.Item(144, 5).Value = 123
$ws.Cells.Item(144, 7).Value = 1
$ws.Cells.Item(144, 8).Value = 3

$ws.Cells.Item(146, 1).Value = "Republica de Chipre"
$ws.Cells.Item(146, 2).Value = 1114
$ws.Cells.Item(146, 3).Value = 24
$ws.Cells.Item(146, 4).Value = 852
$ws.Cells.Item(146, 5).Value = 243
$ws.Cells.Item(146, 8).Value = 19

$ws.Cells.Item(147, 1).Value = "Angola"
$ws.Cells.Item(147, 2).Value = 1109
$ws.Cells.Item(147, 4).Value = 395
$ws.Cells.Item(147, 5).Value = 663
$ws.Cells.Item(147, 8).Value = 51

$ws.Cells.Item(148, 1).Value = "Burkina Faso"
$ws.Cells.Item(148, 2).Value = 1106
$ws.Cells.Item(148, 4).Value = 935
$ws.Cells.Item(148, 5).Value = 118
$ws.Cells.Item(148, 8).Value = 53

$ws.Cells.Item(150, 1).Value = "Togo"
$ws.Cells.Item(150, 2).Value = 927
$ws.Cells.Item(150, 3).Value = 19
$ws.Cells.Item(150, 4).Value = 634
$ws.Cells.Item(150, 5).Value = 275
$ws.Cells.Item(150, 8).Value = 18

$ws.Cells.Item(151, 1).Value = "Principado de Andorra"
$ws.Cells.Item(151, 2).Value = 925
$ws.Cells.Item(151, 3).Value = 3
$ws.Cells.Item(151, 4).Value = 807
$ws.Cells.Item(151, 5).Value = 66
$ws.Cells.Item(151, 8).Value = 52

$ws.Cells.Item(156, 2).Value = 757
$ws.Cells.Item(156, 3).Value = 19
$ws.Cells.Item(156, 4).Value = 237
$ws.Cells.Item(156, 5).Value = 477
$ws.Cells.Item(156, 7).Value = 2
$ws.Cells.Item(156, 8).Value = 43

$ws.Cells.Item(161, 5).Value = 171
$ws.Cells.Item(161, 7).Value = 2
$ws.Cells.Item(161, 8).Value = 2

$ws.Cells.Item(184, 2).Value = 121
$ws.Cells.Item(184, 3).Value = 1
$ws.Cells.Item(184, 5).Value = 13
